$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("B7").Value = '[''Fiat 500 E'', ''Others'']'
$ws.Range("C7").Value = '[0.3, 0.4]'
$ws.Range("D7").Value = '[0.65, 0.7000000000000001]'
$ws.Range("E7").Value = '[8.4, 15.789166666666668]'
$ws.Range("F7").Value = 24.18916666666667

# Row 8
$ws.Range("B8").Value = '[''Opel CORSA'', ''VW E-UP'']'
$ws.Range("C8").Value = '[0.15, 0.2]'
$ws.Range("D8").Value = '[0.9500000000000003, 0.9000000000000002]'
$ws.Range("E8").Value = '[36.000000000000014, 25.760000000000005]'
$ws.Range("F8").Value = 61.76000000000002

# Row 9
$ws.Range("B9").Value = '[''Dacia SPRING'', ''Peugeot 208'', ''Others'']'
$ws.Range("C9").Value = '[0.2, 0.4, 0.4]'
$ws.Range("D9").Value = '[0.9000000000000002, 0.8000000000000002, 1.0]'
$ws.Range("E9").Value = '[18.760000000000005, 18.000000000000007, 31.57833333333333]'
$ws.Range("F9").Value = 68.33833333333334

# Row 11
$ws.Range("B11").Value = '[]'
$ws.Range("C11").Value = '[]'
$ws.Range("D11").Value = '[]'
$ws.Range("E11").Value = '[]'
$ws.Range("F11").Value = 0

# Row 12
$ws.Range("B12").Value = '[''Others'', ''Fiat 500 E'', ''Others'', ''Others'', ''Renault ZOE'']'
$ws.Range("C12").Value = '[0.3, 0.2, 0.25, 0.4, 0.25]'
$ws.Range("D12").Value = '[0.9500000000000003, 0.7500000000000001, 0.65, 1.0, 0.8500000000000002]'
$ws.Range("E12").Value = '[34.20986111111113, 13.200000000000001, 21.052222222222223, 31.57833333333333, 31.20000000000001]'
$ws.Range("F12").Value = 131.2404166666667

# Row 13
$ws.Range("B13").Value = '[''Audi Q4'', ''Others'']'
$ws.Range("C13").Value = '[0.4, 0.45]'
$ws.Range("D13").Value = '[0.7000000000000001, 0.6]'
$ws.Range("E13").Value = '[22.98, 7.8945833333333315]'
$ws.Range("F13").Value = 30.87458333333333

# Row 17
$ws.Range("B17").Value = '[''TESLA MODEL Y'', ''Others'', ''TESLA MODEL Y'']'
$ws.Range("C17").Value = '[0.35, 0.25, 0.2]'
$ws.Range("D17").Value = '[1.0, 0.9000000000000002, 0.7500000000000001]'
$ws.Range("E17").Value = '[48.75, 34.209861111111124, 41.25]'
$ws.Range("F17").Value = 124.2098611111111

# Row 18
$ws.Range("B18").Value = '[''Opel CORSA'', ''Others'']'
$ws.Range("C18").Value = '[0.15, 0.2]'
$ws.Range("D18").Value = '[0.9000000000000002, 0.9000000000000002]'
$ws.Range("E18").Value = '[33.75000000000001, 36.84138888888889]'
$ws.Range("F18").Value = 70.5913888888889

# Row 19
$ws.Range("B19").Value = '[''Others'', ''VW ID.5'']'
$ws.Range("C19").Value = '[0.3, 0.2]'
$ws.Range("D19").Value = '[0.8500000000000002, 0.9500000000000003]'
$ws.Range("E19").Value = '[28.946805555555567, 57.750000000000014]'
$ws.Range("F19").Value = 86.69680555555558

# Row 20
$ws.Range("B20").Value = '[]'
$ws.Range("C20").Value = '[]'
$ws.Range("D20").Value = '[]'
$ws.Range("E20").Value = '[]'
$ws.Range("F20").Value = 0

# Row 31
$ws.Range("B31").Value = '[]'
$ws.Range("C31").Value = '[]'
$ws.Range("D31").Value = '[]'
$ws.Range("E31").Value = '[]'
$ws.Range("F31").Value = 0

# Row 32
$ws.Range("B32").Value = '[''Others'', ''Tesla MODEL 3'', ''SKODA ENYAQ 58kWh'', ''Hyundai IONIQ5 58kWh'']'
$ws.Range("C32").Value = '[0.15, 0.35, 0.25, 0.1]'
$ws.Range("D32").Value = '[0.7500000000000001, 0.9500000000000003, 0.8500000000000002, 0.9000000000000002]'
$ws.Range("E32").Value = '[31.578333333333337, 30.000000000000014, 34.80000000000001, 46.40000000000001]'
$ws.Range("F32").Value = 142.7783333333334

# Row 33
$ws.Range("B33").Value = '[''TESLA MODEL Y'', ''Audi E-TRON'']'
$ws.Range("C33").Value = '[0.05, 0.35]'
$ws.Range("D33").Value = '[0.7000000000000001, 0.8000000000000002]'
$ws.Range("E33").Value = '[48.75, 38.250000000000014]'
$ws.Range("F33").Value = 87.00000000000001

# Row 34
$ws.Range("B34").Value = '[''Others'', ''Others'']'
$ws.Range("C34").Value = '[0.1, 0.3]'
$ws.Range("D34").Value = '[0.8000000000000002, 0.7500000000000001]'
$ws.Range("E34").Value = '[36.84138888888889, 23.683750000000007]'
$ws.Range("F34").Value = 60.5251388888889

# Row 36
$ws.Range("B36").Value = '[''SKODA ENYAQ 58kWh'', ''SKODA ENYAQ 77kWh'', ''Fiat 500 E'']'
$ws.Range("C36").Value = '[0.2, 0.25, 0.3]'
$ws.Range("D36").Value = '[0.9000000000000002, 1.0, 0.9500000000000003]'
$ws.Range("E36").Value = '[40.60000000000001, 57.75, 15.600000000000009]'
$ws.Range("F36").Value = 113.95

# Row 37
$ws.Range("B37").Value = '[''Fiat 500 E'', ''Fiat 500 E'', ''Renault ZOE'']'
$ws.Range("C37").Value = '[0.2, 0.3, 0.2]'
$ws.Range("D37").Value = '[0.8500000000000002, 0.8000000000000002, 0.9000000000000002]'
$ws.Range("E37").Value = '[15.600000000000003, 12.000000000000005, 36.400000000000006]'
$ws.Range("F37").Value = 64.00000000000001

# Row 39
$ws.Range("B39").Value = '[]'
$ws.Range("C39").Value = '[]'
$ws.Range("D39").Value = '[]'
$ws.Range("E39").Value = '[]'
$ws.Range("F39").Value = 0

# Row 41
$ws.Range("B41").Value = '[''Hyundai KONA 39 kWh'', ''MINI Cooper SE'']'
$ws.Range("C41").Value = '[0.05, 0.2]'
$ws.Range("D41").Value = '[0.9500000000000003, 0.9500000000000003]'
$ws.Range("E41").Value = '[35.10000000000001, 21.675000000000004]'
$ws.Range("F41").Value = 56.77500000000001

# Row 42
$ws.Range("B42").Value = '[]'
$ws.Range("C42").Value = '[]'
$ws.Range("D42").Value = '[]'
$ws.Range("E42").Value = '[]'
$ws.Range("F42").Value = 0

# Row 43
$ws.Range("B43").Value = '[''TESLA MODEL Y'', ''Others'', ''Others'', ''Tesla MODEL 3'']'
$ws.Range("C43").Value = '[0.05, 0.25, 0.05, 0.05]'
$ws.Range("D43").Value = '[0.7500000000000001, 0.9500000000000003, 0.7000000000000001, 1.0]'
$ws.Range("E43").Value = '[52.50000000000001, 36.8413888888889, 34.20986111111111, 47.5]'
$ws.Range("F43").Value = 171.05125

# Row 44
$ws.Range("B44").Value = '[''Fiat 500 E'']'
$ws.Range("C44").Value = '[0.35]'
$ws.Range("D44").Value = '[0.9000000000000002]'
$ws.Range("E44").Value = '[13.200000000000006]'
$ws.Range("F44").Value = 13.20000000000001
